# The ER (entity-relationship) mapping table for "ferias" had its row for
# the "status" field pointing to a placeholder column name. Update the
# mapped (new) column name for that row from "request_vacations" to "XXX".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "XXX"
